# Update gh-pages output data across sheets

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 1410
$wsExpo.Range("F3").Value = 2963
$wsExpo.Range("F4").Value = 26
$wsExpo.Range("F5").Value = 273

# Sheet "演出" (performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("G2").Value = "不可售"

# Sheet "全部类型" (all types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("G2").Value = "不可售"
$wsAll.Range("F3").Value = 1410
$wsAll.Range("F4").Value = 2963
$wsAll.Range("F5").Value = 26
$wsAll.Range("F7").Value = 273
